$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update Price (D) and Volume/1h (E) columns.
# Some Price values are plain decimals (e.g. "313.84"); force those cells to
# Text format first so Excel keeps them as strings instead of parsing them
# as numbers, matching the source data feed formatting.
$ws.Range("D2").Value = "27.249.98"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.856.27"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.84"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07302"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8911"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +2.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07878"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "1.853.75"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.408"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.512"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.75"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "27.287.60"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.077"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "2.116.80"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.031"
$ws.Range("E25").Value = "  +9.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.75"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.99"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.050"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08825"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.150"
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7705"
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.525"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.721"
$ws.Range("E36").Value = "  +10.66%  "
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01942"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05230"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.075"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5119"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.550"
$ws.Range("E44").Value = "  +5.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4786"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.35"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.94"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.647"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06206"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.65"
$ws.Range("E51").Value = "  +2.09%  "
